$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update SP2D GUP 17 realisasi/sisa figures (existing rows) ---
$ws.Range("N7").Value = 475000000
$ws.Range("O7").Value = 95000000

$ws.Range("N50").Value = 80856500
$ws.Range("O50").Value = 74307500

$ws.Range("N58").Value = 31500000
$ws.Range("O58").Value = 0

$ws.Range("N59").Value = 30000000
$ws.Range("O59").Value = 0

# --- Fill in newly reported contract rows 67-69 ---
# Date-like text columns (E, F, K, L) must be entered with a leading
# apostrophe to keep Excel from auto-converting them to real dates; an
# extra leading apostrophe is used when the target text itself should
# retain a literal apostrophe character. ClearFormats() afterwards drops
# the quote-prefix formatting so the cells end up as plain shared-string
# text cells (matching the rest of the sheet).

# Row 67
$ws.Range("A67").Value = 63
$ws.Range("B67").Value = 626402
$ws.Range("C67").Value = "A/175.22008802/0/0"
$ws.Range("D67").Value = "PT. BIGJEK APLIKASI MANDIRI"
$ws.Range("E67").Value = "''07-NOV-22"
$ws.Range("F67").Value = "''04-NOV-22"
$ws.Range("G67").Value = "Tidak terlambat"
$ws.Range("H67").Value = "SPK-4211/PPK/BRSDM.1/XI/2022"
$ws.Range("I67").Value = "Pekerjaan Jasa Pengembangan Dashboard Data Center BRSDM"
$ws.Range("J67").Value = 522191
$ws.Range("K67").Value = "''04-NOV-22"
$ws.Range("L67").Value = "'18-NOV-22"
$ws.Range("E67:L67").ClearFormats()
$ws.Range("M67").Value = 99900000
$ws.Range("N67").Value = 0
$ws.Range("O67").Value = 99900000

# Row 68
$ws.Range("A68").Value = 64
$ws.Range("B68").Value = 626402
$ws.Range("C68").Value = "A/175.22008962/0/0"
$ws.Range("D68").Value = "PT. GITA AGUNG PRATAMA"
$ws.Range("E68").Value = "'10-NOV-22"
$ws.Range("F68").Value = "''08-NOV-22"
$ws.Range("G68").Value = "Tidak terlambat"
$ws.Range("H68").Value = "SPK-3752/PPK/BRSDM.1/XI/2022"
$ws.Range("I68").Value = "Pengadaan Mebelair Ruang Kepala BRSDM"
$ws.Range("J68").Value = 532111
$ws.Range("K68").Value = "''08-NOV-22"
$ws.Range("L68").Value = "'16-NOV-22"
$ws.Range("E68:L68").ClearFormats()
$ws.Range("M68").Value = 191440000
$ws.Range("N68").Value = 0
$ws.Range("O68").Value = 191440000

# Row 69
$ws.Range("A69").Value = 65
$ws.Range("B69").Value = 626402
$ws.Range("C69").Value = "A/175.22008963/0/0"
$ws.Range("D69").Value = "CV. GRAHA IDEKU INTERIOR"
$ws.Range("E69").Value = "'10-NOV-22"
$ws.Range("F69").Value = "''08-NOV-22"
$ws.Range("G69").Value = "Tidak terlambat"
$ws.Range("H69").Value = "SPK-3673/PPK/BRSDM.1/XI/2022"
$ws.Range("I69").Value = "Pekerjaan Pembuatan dan Pemasangan Funiture Joinery Ruang Sekretaris Kepala BRSDM"
$ws.Range("J69").Value = 533121
$ws.Range("K69").Value = "''08-NOV-22"
$ws.Range("L69").Value = "'16-NOV-22"
$ws.Range("E69:L69").ClearFormats()
$ws.Range("M69").Value = 146381000
$ws.Range("N69").Value = 0
$ws.Range("O69").Value = 146381000
